$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '23.238.25'
$ws.Range('E2').Value = '  +0.80%  '

# Row 3
$ws.Range('D3').Value = '1.603.07'
$ws.Range('E3').Value = '  +0.07%  '

# Row 4
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('E5').Value = '  -0.16%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '303.42'
$ws.Range('E6').Value = '  +0.78%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3776'
$ws.Range('E7').Value = '  -0.03%  '

# Row 8
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '51.77'
$ws.Range('E8').Value = '  +3.59%  '

# Row 9
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3639'
$ws.Range('E9').Value = '  -0.18%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.273'
$ws.Range('E10').Value = '  +0.92%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08138'
$ws.Range('E11').Value = '  +0.15%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.22%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.87'
$ws.Range('E13').Value = '  +0.35%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.606'
$ws.Range('E14').Value = '  +0.29%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.426'
$ws.Range('E15').Value = '  +0.70%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001250'
$ws.Range('E16').Value = '  -0.34%  '

# Row 17
$ws.Range('D17').Value = '1.606.91'
$ws.Range('E17').Value = '  +0.52%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.00'
$ws.Range('E18').Value = '  +2.05%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06901'
$ws.Range('E19').Value = '  +0.18%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.17'
$ws.Range('E20').Value = '  -0.53%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.535'
$ws.Range('E21').Value = '  -0.48%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9998'

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.98'
$ws.Range('E23').Value = '  -0.98%  '

# Row 24
$ws.Range('D24').Value = '23.231.41'
$ws.Range('E24').Value = '  +0.72%  '

# Row 25
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.015'
$ws.Range('E25').Value = '  +8.11%  '

# Row 26
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.391'
$ws.Range('E26').Value = '  +0.84%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '21.23'
$ws.Range('E27').Value = '  +0.71%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '150.25'
$ws.Range('E28').Value = '  -0.09%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.256'
$ws.Range('E29').Value = '  -0.06%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '134.24'
$ws.Range('E30').Value = '  +0.31%  '

# Row 31
$ws.Range('E31').Value = '  +1.06%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.779'
$ws.Range('E32').Value = '  -0.98%  '

# Row 33
$ws.Range('D33').Value = '1.781.91'
$ws.Range('E33').Value = '  +0.35%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9662'
$ws.Range('E34').Value = '  +0.64%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07547'
$ws.Range('E35').Value = '  -1.64%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02741'
$ws.Range('E36').Value = '  +0.99%  '

# Row 37
$ws.Range('E37').Value = '  -2.24%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2536'
$ws.Range('E38').Value = '  -0.16%  '

# Row 39
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.117'
$ws.Range('E39').Value = '  -2.87%  '

# Row 40
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.08812'
$ws.Range('E40').Value = '  -1.08%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.392'
$ws.Range('E41').Value = '  +1.76%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.7134'
$ws.Range('E42').Value = '  +1.15%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '12.57'
$ws.Range('E43').Value = '  -0.48%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.66'
$ws.Range('E44').Value = '  +2.15%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6549'
$ws.Range('E45').Value = '  -1.10%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.329'
$ws.Range('E46').Value = '  +0.12%  '

# Row 47
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.9992'
$ws.Range('E47').Value = '  -0.13%  '

# Row 48
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.015'
$ws.Range('E48').Value = '  +0.43%  '

# Row 49
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '132.69'
$ws.Range('E49').Value = '  +0.14%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07963'
$ws.Range('E50').Value = '  +0.45%  '

# Row 51
$ws.Range('B51').Value = 'Flow'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('E51').Value = '  -2.69%  '
